$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 583, pushing existing rows 583-622 down to 587-626
$ws.Rows("583:586").Insert()

# Common (fixed) values shared by every row in this block
$mercado   = 2
$region    = "Comercializadora del Agro de Limarí"
$codreg    = "Coquimbo"
$categoriaId = 4
$varId     = 100112043
$categoria = "Pepino dulce"
$variedad  = "Cultivar IV Región"
$unidad    = "$/bandeja 18 kilos"
$origen    = "Provincia de Limarí"
$kgUnid    = 18
$clasif    = "Hortaliza"
$fecha     = 45021

# Row 583 - Especial
$ws.Cells.Item(583, 1).Value = $mercado
$ws.Cells.Item(583, 2).Value = $region
$ws.Cells.Item(583, 3).Value = $codreg
$ws.Cells.Item(583, 4).Value = $fecha
$ws.Cells.Item(583, 5).Value = $categoriaId
$ws.Cells.Item(583, 6).Value = $varId
$ws.Cells.Item(583, 7).Value = $categoria
$ws.Cells.Item(583, 8).Value = $variedad
$ws.Cells.Item(583, 9).Value = "Especial"
$ws.Cells.Item(583, 10).Value = 600
$ws.Cells.Item(583, 11).Value = 11000
$ws.Cells.Item(583, 12).Value = 12000
$ws.Cells.Item(583, 13).Value = 11500
$ws.Cells.Item(583, 14).Value = $unidad
$ws.Cells.Item(583, 15).Value = $origen
$ws.Cells.Item(583, 16).Value = 639
$ws.Cells.Item(583, 17).Value = $kgUnid
$ws.Cells.Item(583, 18).Value = $clasif

# Row 584 - Primera
$ws.Cells.Item(584, 1).Value = $mercado
$ws.Cells.Item(584, 2).Value = $region
$ws.Cells.Item(584, 3).Value = $codreg
$ws.Cells.Item(584, 4).Value = $fecha
$ws.Cells.Item(584, 5).Value = $categoriaId
$ws.Cells.Item(584, 6).Value = $varId
$ws.Cells.Item(584, 7).Value = $categoria
$ws.Cells.Item(584, 8).Value = $variedad
$ws.Cells.Item(584, 9).Value = "Primera"
$ws.Cells.Item(584, 10).Value = 800
$ws.Cells.Item(584, 11).Value = 9000
$ws.Cells.Item(584, 12).Value = 10000
$ws.Cells.Item(584, 13).Value = 9500
$ws.Cells.Item(584, 14).Value = $unidad
$ws.Cells.Item(584, 15).Value = $origen
$ws.Cells.Item(584, 16).Value = 528
$ws.Cells.Item(584, 17).Value = $kgUnid
$ws.Cells.Item(584, 18).Value = $clasif

# Row 585 - Segunda
$ws.Cells.Item(585, 1).Value = $mercado
$ws.Cells.Item(585, 2).Value = $region
$ws.Cells.Item(585, 3).Value = $codreg
$ws.Cells.Item(585, 4).Value = $fecha
$ws.Cells.Item(585, 5).Value = $categoriaId
$ws.Cells.Item(585, 6).Value = $varId
$ws.Cells.Item(585, 7).Value = $categoria
$ws.Cells.Item(585, 8).Value = $variedad
$ws.Cells.Item(585, 9).Value = "Segunda"
$ws.Cells.Item(585, 10).Value = 700
$ws.Cells.Item(585, 11).Value = 7000
$ws.Cells.Item(585, 12).Value = 8000
$ws.Cells.Item(585, 13).Value = 7500
$ws.Cells.Item(585, 14).Value = $unidad
$ws.Cells.Item(585, 15).Value = $origen
$ws.Cells.Item(585, 16).Value = 417
$ws.Cells.Item(585, 17).Value = $kgUnid
$ws.Cells.Item(585, 18).Value = $clasif

# Row 586 - Tercera
$ws.Cells.Item(586, 1).Value = $mercado
$ws.Cells.Item(586, 2).Value = $region
$ws.Cells.Item(586, 3).Value = $codreg
$ws.Cells.Item(586, 4).Value = $fecha
$ws.Cells.Item(586, 5).Value = $categoriaId
$ws.Cells.Item(586, 6).Value = $varId
$ws.Cells.Item(586, 7).Value = $categoria
$ws.Cells.Item(586, 8).Value = $variedad
$ws.Cells.Item(586, 9).Value = "Tercera"
$ws.Cells.Item(586, 10).Value = 500
$ws.Cells.Item(586, 11).Value = 5000
$ws.Cells.Item(586, 12).Value = 6000
$ws.Cells.Item(586, 13).Value = 5500
$ws.Cells.Item(586, 14).Value = $unidad
$ws.Cells.Item(586, 15).Value = $origen
$ws.Cells.Item(586, 16).Value = 306
$ws.Cells.Item(586, 17).Value = $kgUnid
$ws.Cells.Item(586, 18).Value = $clasif

# Apply the same date style (numFmt) as the rest of column D to the newly inserted date cells
$ws.Range("D583:D586").NumberFormat = $ws.Range("D587").NumberFormat
